{"js": "// Change 1: the deployed WAR file was renamed from \"specchio-services.war\"\n// to \"specchio-webapp.war\" (only the middle word actually changed).\nconst body = context.document.body;\n\nconst warResults = body.search(\"specchio-services.war\", { matchCase: true });\nwarResults.load(\"items\");\nawait context.sync();\n\nif (warResults.items.length > 0) {\n  warResults.items[0].insertText(\"specchio-webapp.war\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Change 2: Word re-anchors its hidden \"_GoBack\" bookmark (last-edit marker)\n// at the location of the most recent edit, which is immediately after the\n// sentence \"...usually specchio_service.\" Remove the old one (it used to\n// sit in the empty paragraph at the very end of the document) and insert a\n// fresh, empty bookmark right after that sentence.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst goBackAnchor = body.search(\"usually specchio_service.\", { matchCase: true });\ngoBackAnchor.load(\"items\");\nawait context.sync();\n\nif (goBackAnchor.items.length > 0) {\n  goBackAnchor.items[0].getRange(Word.RangeLocation.end).insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// Change 3: the \"Upgrading from V3.1.0 to V3.2.0\" heading used to be split\n// across two runs (\"...to V3.2\" + \".0\"); re-write it as a single run with\n// the same resulting text so the runs collapse back into one.\nconst headingResults = body.search(\"Upgrading from V3.1.0 to V3.2.0\", { matchCase: true });\nheadingResults.load(\"items\");\nawait context.sync();\n\nif (headingResults.items.length > 0) {\n  headingResults.items[0].insertText(\"Upgrading from V3.1.0 to V3.2.0\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Change 1: the deployed WAR file was renamed from \"specchio-services.war\"\n# to \"specchio-webapp.war\" (only the middle word actually changed).\n$warRange = $d.Content\n$warRange.Find.Text = \"specchio-services.war\"\n$warRange.Find.Execute() | Out-Null\nif ($warRange.Find.Found) {\n    $warRange.Text = \"specchio-webapp.war\"\n}\n\n# Change 2: Word re-anchors its hidden \"_GoBack\" bookmark (last-edit marker)\n# at the location of the most recent edit, which is immediately after the\n# sentence \"...usually specchio_service.\" Remove the old one (it used to\n# sit in the empty paragraph at the very end of the document) and insert a\n# fresh, empty bookmark right after that sentence.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$anchor = $d.Content\n$anchor.Find.Text = \"usually specchio_service.\"\n$anchor.Find.Execute() | Out-Null\nif ($anchor.Find.Found) {\n    $anchor.Collapse(0)  # wdCollapseEnd - collapse to the end of the match\n\n    # A genuinely zero-length range can't host a new bookmark directly, so\n    # insert a one-character placeholder, bookmark that character, then\n    # delete the placeholder text - this leaves the (now empty) bookmark\n    # correctly anchored between the \".\" run and the end of the paragraph.\n    $anchor.InsertAfter(\"X\")\n    $placeholder = $d.Range($anchor.Start, $anchor.Start + 1)\n    $placeholder.Bookmarks.Add(\"_GoBack\")\n    $goBack = $d.Bookmarks(\"_GoBack\")\n    $goBack.Range.Text = \"\"\n}\n\n# Change 3: the \"Upgrading from V3.1.0 to V3.2.0\" heading used to be split\n# across two runs (\"...to V3.2\" + \".0\"); re-write it as a single run with\n# the same resulting text so the runs collapse back into one.\n$heading = $d.Content\n$heading.Find.Text = \"Upgrading from V3.1.0 to V3.2.0\"\n$heading.Find.Execute() | Out-Null\nif ($heading.Find.Found) {\n    $heading.Text = \"\"\n    $heading.InsertAfter(\"Upgrading from V3.1.0 to V3.2.0\")\n}\n"}
